# feat: remove unused analysis
#
# Drops the "Public/Private/ONG" breakdown rows (original rows 8-10 and
# 17-19, i.e. Lethality_Public_Tx, Lethality_Private_Tx, Lethality_ONG_Tx,
# Lethality_uti_Public_Tx, Lethality_uti_Private_Tx, Lethality_uti_ONG_Tx)
# from the lethality-rates table. The remaining analyses (uti_* and
# non_uti_*) are moved up into rows 8-19, while keeping each row's original
# index (column A) and row number untouched, and the now-empty trailing
# rows 20-25 are removed. The sheet shrinks from A1:I25 to A1:I19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = "Lethality_uti_tx"
$ws.Range("C8").Value = -0.3431171880430295
$ws.Range("D8").Value = -2.142180710707942
$ws.Range("E8").Value = 1.489021152578096
$ws.Range("F8").Value = 39867
$ws.Range("G8").Value = 39901
$ws.Range("H8").Value = 39592
$ws.Range("I8").Value = -0.6897935635989666
$ws.Range("B9").Value = "Lethality_uti_S_Tx"
$ws.Range("C9").Value = 8.70574001467217
$ws.Range("D9").Value = 4.187997486319572
$ws.Range("E9").Value = 13.41937840478338
$ws.Range("F9").Value = 34092
$ws.Range("G9").Value = 34053
$ws.Range("H9").Value = 40388
$ws.Range("I9").Value = 18.46767570104423
$ws.Range("B10").Value = "Lethality_uti_N_Tx"
$ws.Range("C10").Value = -1.352338033009881
$ws.Range("D10").Value = -8.693756849846334
$ws.Range("E10").Value = 6.579362766577401
$ws.Range("F10").Value = 51552
$ws.Range("G10").Value = 56526
$ws.Range("H10").Value = 49989
$ws.Range("I10").Value = -3.031890130353818
$ws.Range("B11").Value = "Lethality_uti_NE_Tx"
$ws.Range("C11").Value = -1.352338033009881
$ws.Range("D11").Value = -8.693756849846334
$ws.Range("E11").Value = 6.579362766577401
$ws.Range("F11").Value = 51552
$ws.Range("G11").Value = 56526
$ws.Range("H11").Value = 49989
$ws.Range("I11").Value = -3.031890130353818
$ws.Range("B12").Value = "Lethality_uti_SE_Tx"
$ws.Range("C12").Value = 0.09068439544728246
$ws.Range("D12").Value = -2.574103687616813
$ws.Range("E12").Value = 2.828359624500498
$ws.Range("F12").Value = 37389
$ws.Range("G12").Value = 36724
$ws.Range("H12").Value = 37431
$ws.Range("I12").Value = 0.1123325042124689
$ws.Range("B13").Value = "Lethality_uti_CO_Tx"
$ws.Range("C13").Value = -3.21181878395681
$ws.Range("D13").Value = -8.565486054774496
$ws.Range("E13").Value = 2.455316038772337
$ws.Range("F13").Value = 42259
$ws.Range("G13").Value = 42918
$ws.Range("H13").Value = 39330
$ws.Range("I13").Value = -6.931067938190681
$ws.Range("B14").Value = "Lethality_non_uti_Tx"
$ws.Range("C14").Value = 14.9414357382152
$ws.Range("D14").Value = 13.54164374982032
$ws.Range("E14").Value = 16.35848498611474
$ws.Range("F14").Value = 6535
$ws.Range("G14").Value = 7843
$ws.Range("H14").Value = 8607
$ws.Range("I14").Value = 31.7061973986228
$ws.Range("B15").Value = "Lethality_non_uti_S_Tx"
$ws.Range("C15").Value = 30.23212302935951
$ws.Range("D15").Value = 26.39844905537223
$ws.Range("E15").Value = 34.18207260837727
$ws.Range("F15").Value = 4698
$ws.Range("G15").Value = 6038
$ws.Range("H15").Value = 7973
$ws.Range("I15").Value = 69.71051511281397
$ws.Range("B16").Value = "Lethality_non_uti_N_Tx"
$ws.Range("C16").Value = 9.843325733999686
$ws.Range("D16").Value = 4.750020257837484
$ws.Range("E16").Value = 15.18428520211006
$ws.Range("F16").Value = 5915
$ws.Range("G16").Value = 6841
$ws.Range("H16").Value = 7104
$ws.Range("I16").Value = 20.10143702451395
$ws.Range("B17").Value = "Lethality_non_uti_NE_Tx"
$ws.Range("C17").Value = 9.843325733999686
$ws.Range("D17").Value = 4.750020257837484
$ws.Range("E17").Value = 15.18428520211006
$ws.Range("F17").Value = 5915
$ws.Range("G17").Value = 6841
$ws.Range("H17").Value = 7104
$ws.Range("I17").Value = 20.10143702451395
$ws.Range("B18").Value = "Lethality_non_uti_SE_Tx"
$ws.Range("C18").Value = 13.53710191105546
$ws.Range("D18").Value = 11.56836256871352
$ws.Range("E18").Value = 15.54058169870687
$ws.Range("F18").Value = 8322
$ws.Range("G18").Value = 9707
$ws.Range("H18").Value = 10713
$ws.Range("I18").Value = 28.73107426099495
$ws.Range("B19").Value = "Lethality_non_uti_CO_Tx"
$ws.Range("C19").Value = 17.65726333918107
$ws.Range("D19").Value = 11.64161622293005
$ws.Range("E19").Value = 23.99705490489075
$ws.Range("F19").Value = 5267
$ws.Range("G19").Value = 5364
$ws.Range("H19").Value = 7405
$ws.Range("I19").Value = 40.59236757167267

# Remove the now-duplicated trailing rows (originally rows 20-25)
$ws.Range("A20:I25").Delete()
